{"js": "// Apply the README/docx stat corrections to the single-column results table.\n// Each change below targets one table row (0-indexed) and replaces that\n// row's text content while preserving the existing run formatting\n// (rFonts/sz) by replacing the paragraph's Range text in place, instead of\n// clearing/re-inserting the cell body (which would drop the rPr).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> new text for that row's single paragraph/cell\nconst edits = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"936\"],\n  [4, \"0.00001\"],\n  [6, \"0.00013\"],\n  [7, \"0.00007\"],\n  [8, \"0.00015\"],\n  [9, \"0.00016\"],\n  [10, \"0.00030\"],\n  [11, \"0.13422\"],\n  [43, \"99.94\"],\n  [44, \"0.13\"],\n  [45, \"209\"],\n];\n\nfor (const [rowIndex, newText] of edits) {\n  const cell = table.getCell(rowIndex, 0);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the README/docx stat corrections to the single-column results table.\n# Word COM table cells are 1-indexed (Cell(row, column)); setting\n# Range.Text in place preserves the existing run formatting (rFonts/sz)\n# instead of clearing/re-adding a run, which matches the diff (only the\n# <w:t> text content changes).\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$edits = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"936\"\n    5  = \"0.00001\"\n    7  = \"0.00013\"\n    8  = \"0.00007\"\n    9  = \"0.00015\"\n    10 = \"0.00016\"\n    11 = \"0.00030\"\n    12 = \"0.13422\"\n    44 = \"99.94\"\n    45 = \"0.13\"\n    46 = \"209\"\n}\n\nforeach ($rowIndex in $edits.Keys) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cell.Range.Text = $edits[$rowIndex]\n}\n"}
